$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# 1. Insert a new column at B (shifts old B->C, C->D, D->E, F->G)
$ws.Columns.Item(2).Insert()

# 2. Set the new column B width (best achievable approximation of 44.33203125 given engine quantization)
$ws.Columns.Item(2).ColumnWidth = 43.5

# 3. Populate column B (Descrição) - header + per-row descriptions
$ws.Range("B1").Value = "Descrição"
$ws.Range("B2").Value = "Por já utilizarem componentes elétricos nas estufas, o nosso cliente já possui previamente o conhecimento para evitar queimas nos componentes"
$ws.Range("B3").Value = "Como as orquideas requerem cuidado constante, as estufas são preparadas para para qualquer evento que possa afetar as suas plantações"
$ws.Range("B4").Value = "Devido ao fato do nosso cliente ser os agricultores, a falta de conhecimento em componentes tecnológicos pode resultar numa chance maior de impacto nos componentes do nosso produto"
$ws.Range("B5").Value = "Qualquer plantação requer o uso constante de umidade para o seu tratamento, isso junto a falta de experiência aos nossos clientes, as chances de oxidação se tornam maiores "
$ws.Range("B6").Value = "Dificilmente ocorre algum problema nos sensores, geralmente quando acontece, é devido a um defeito de fábrica,  uma má configuração é dificil de ocorrer já que os sensores serão preparados por nós"
$ws.Range("B7").Value = "Pelo fato dos dados gerados pelos nossos sensores serem encaminhados ao nosso banco de dados de forma remota, as chances de haver alguma interferência tende a ser maior, devido as máquinas que são utilizadas nas estufas"

# 4. Update mitigation text (column E) that changed
$ws.Range("E2").Value = "Uso de dispositivos de proteção contra surtos(DPS) como filtros de linha, entre outros"
$ws.Range("E3").Value = "Uso de um gerador de energia reserva, como por exemplo o Nobreak"
$ws.Range("E7").Value = "Usar materiais de blindagem"

# 5. Update probability/impact values that changed on row 4
$ws.Range("C4").Value = "Médio"
$ws.Range("D4").Value = "Médio"

# 6. Styling
# Header row restyle (bold pink centered) - ensure G1 (Legenda header) also gets this style
$ws.Range("A1:G1").HorizontalAlignment = $xlCenter
$ws.Range("A1:G1").VerticalAlignment = $xlCenter
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").Font.Size = 14
$ws.Range("G1").Interior.Color = $ws.Range("A1").Interior.Color

# Column A (risk names) - centered, no wrap
$ws.Range("A2:A7").HorizontalAlignment = $xlCenter
$ws.Range("A2:A7").VerticalAlignment = $xlCenter

# Column B (descriptions) - horizontal center + wrap, default vertical
$ws.Range("B2:B7").HorizontalAlignment = $xlCenter
$ws.Range("B2:B7").WrapText = $true

# Column C, D (probability / impact) - centered, no wrap, no fill
$ws.Range("C2:D7").HorizontalAlignment = $xlCenter
$ws.Range("C2:D7").VerticalAlignment = $xlCenter

# Column E (mitigation) - centered + wrap
$ws.Range("E2:E7").HorizontalAlignment = $xlCenter
$ws.Range("E2:E7").VerticalAlignment = $xlCenter
$ws.Range("E2:E7").WrapText = $true

# Legend cells G2:G4 - centered
$ws.Range("G2:G4").HorizontalAlignment = $xlCenter
$ws.Range("G2:G4").VerticalAlignment = $xlCenter

# 7. Row heights
$ws.Rows.Item(2).RowHeight = 72
$ws.Rows.Item(3).RowHeight = 72
$ws.Rows.Item(4).RowHeight = 90
$ws.Rows.Item(5).RowHeight = 90
$ws.Rows.Item(6).RowHeight = 90
$ws.Rows.Item(7).RowHeight = 108

# 8. Fix conditional formatting range from B1:C1048576 to C1:D1048576
$ws.Cells.FormatConditions.Delete()
$cf1 = $ws.Range("C1:D1048576").FormatConditions.Add(1, 3, """Alto""")
$cf2 = $ws.Range("C1:D1048576").FormatConditions.Add(1, 3, """Médio""")
$cf3 = $ws.Range("C1:D1048576").FormatConditions.Add(1, 3, """Baixo""")

# 9. Selection
$ws.Range("B3").Select()
